$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = @'
Pipeline(steps=[('scaler', StandardScaler()), ('selector', None),
                ('model',
                 SVC(C=3, class_weight='balanced', kernel='sigmoid',
                     random_state=42))])
'@
$ws.Range("B2").Value = 0.6952380952380952
$ws.Range("C2").Value = @'
{'selector': None, 'scaler': StandardScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': 'balanced', 'model__C': 3}
'@
$ws.Range("D2").Value = 0.2222222222222222
$ws.Range("E2").Value = @'
[1 1 0 0 1 0 0 0 0 1 0 1]
'@
$ws.Range("F2").Value = @'
[1 0 0 0 0 1 1 0 1 0 0 0]
'@
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.6776371308016877
$ws.Range("I2").Value = 0.03490433830610162
$ws.Range("J2").Value = 0.5885473176612416
$ws.Range("K2").Value = 0.06210767213185699

# Row 3
$ws.Range("A3").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a642ef7c0>),
                ('model',
                 SVC(C=0.0001, class_weight='balanced', kernel='linear',
                     random_state=42))])
'@
$ws.Range("B3").Value = 0.619047619047619
$ws.Range("C3").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6414dbb0>, 'scaler': None, 'model__kernel': 'linear', 'model__class_weight': 'balanced', 'model__C': 0.0001}
'@
$ws.Range("D3").Value = 0.7777777777777778
$ws.Range("E3").Value = @'
[1 1 0 1 0 0 1 0 1 1 1 0]
'@
$ws.Range("F3").Value = @'
[1 1 1 1 1 1 1 1 1 1 1 0]
'@
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.6609465737514517
$ws.Range("I3").Value = 0.02026234985050712
$ws.Range("J3").Value = 0.5414634146341463
$ws.Range("K3").Value = 0.05748756504422877

# Row 4
$ws.Range("A4").Value = @'
Pipeline(steps=[('scaler', None), ('selector', None),
                ('model',
                 SVC(C=5, class_weight='balanced', kernel='sigmoid',
                     random_state=42))])
'@
$ws.Range("B4").Value = 0.6285714285714286
$ws.Range("C4").Value = @'
{'selector': None, 'scaler': None, 'model__kernel': 'sigmoid', 'model__class_weight': 'balanced', 'model__C': 5}
'@
$ws.Range("D4").Value = 0.4285714285714285
$ws.Range("E4").Value = @'
[1 0 1 1 1 1 0 1 0 1 0 1]
'@
$ws.Range("F4").Value = @'
[0 1 1 0 1 0 1 1 1 0 0 0]
'@
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6214888487040385
$ws.Range("I4").Value = 0.02568081191328808
$ws.Range("J4").Value = 0.5398432790837855
$ws.Range("K4").Value = 0.05926201853284085

# Row 5
$ws.Range("A5").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model', SVC(C=5, kernel='sigmoid', random_state=42))])
'@
$ws.Range("B5").Value = 0.6857142857142857
$ws.Range("C5").Value = @'
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 5}
'@
$ws.Range("D5").Value = 0.5333333333333333
$ws.Range("E5").Value = @'
[1 1 0 0 0 0 1 0 1 1 1 1]
'@
$ws.Range("F5").Value = @'
[0 0 1 0 1 1 1 1 1 0 1 1]
'@
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.6422757475083056
$ws.Range("I5").Value = 0.02675775831887378
$ws.Range("J5").Value = 0.5428571428571429
$ws.Range("K5").Value = 0.05537779399289352

# Row 6
$ws.Range("A6").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 SVC(C=0.0001, class_weight='balanced', kernel='linear',
                     random_state=42))])
'@
$ws.Range("B6").Value = 0.6761904761904761
$ws.Range("C6").Value = @'
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__kernel': 'linear', 'model__class_weight': 'balanced', 'model__C': 0.0001}
'@
$ws.Range("D6").Value = 0.4615384615384615
$ws.Range("E6").Value = @'
[1 1 1 1 0 0 0 0 1 1 0 0]
'@
$ws.Range("F6").Value = @'
[1 1 0 0 0 1 1 1 1 0 1 0]
'@
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.6605758582502769
$ws.Range("I6").Value = 0.03745735103361374
$ws.Range("J6").Value = 0.5551495016611295
$ws.Range("K6").Value = 0.05991010253942086

# Reset auto row-height side effects from multi-line values so rows keep
# the workbook's default (no explicit customHeight), matching the source.
for ($r = 2; $r -le 6; $r++) {
    $ws.Rows.Item($r).EntireRow.AutoFit()
}
